$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1275
$ws.Cells.Item(2, 9).Value = 1350.5
$ws.Cells.Item(2, 10).Value = 1237.25
$ws.Cells.Item(2, 11).Value = 1350.5
$ws.Cells.Item(2, 12).Value = 1237.25
$ws.Cells.Item(2, 13).Value = -1237.5
$ws.Cells.Item(2, 14).Value = -1463.25
$ws.Cells.Item(3, 8).Value = 39933
$ws.Cells.Item(3, 10).Value = 39933
$ws.Cells.Item(3, 12).Value = 39933
$ws.Cells.Item(3, 14).Value = -40161
$ws.Cells.Item(17, 8).Value = 1311.9584
$ws.Cells.Item(17, 10).Value = 1311.9584
$ws.Cells.Item(17, 12).Value = 3935.8752
$ws.Cells.Item(17, 14).Value = -4271.8752
$ws.Cells.Item(19, 8).Value = 2380.9524
$ws.Cells.Item(19, 9).Value = 2102.7693
$ws.Cells.Item(19, 10).Value = 2833
$ws.Cells.Item(19, 11).Value = 2102.7693
$ws.Cells.Item(19, 12).Value = 2833
$ws.Cells.Item(19, 13).Value = -1927.7693
$ws.Cells.Item(19, 14).Value = -3183
$ws.Cells.Item(21, 8).Value = 8375
$ws.Cells.Item(21, 9).Value = 8375
$ws.Cells.Item(21, 11).Value = 8375
$ws.Cells.Item(21, 13).Value = -7907
$ws.Cells.Item(23, 8).Value = 8375
$ws.Cells.Item(23, 9).Value = 8375
$ws.Cells.Item(23, 11).Value = 8375
$ws.Cells.Item(23, 13).Value = -8141
$ws.Cells.Item(29, 8).Value = 7699.25
$ws.Cells.Item(29, 9).Value = 4399.5
$ws.Cells.Item(29, 10).Value = 10999
$ws.Cells.Item(29, 11).Value = 13198.5
$ws.Cells.Item(29, 12).Value = 32997
$ws.Cells.Item(29, 13).Value = -12917.5
$ws.Cells.Item(29, 14).Value = -33559
$ws.Cells.Item(32, 8).Value = 13227.615
$ws.Cells.Item(32, 9).Value = 13896.2
$ws.Cells.Item(32, 10).Value = 10999
$ws.Cells.Item(32, 11).Value = 13896.2
$ws.Cells.Item(32, 12).Value = 10999
$ws.Cells.Item(32, 13).Value = -13570.2
$ws.Cells.Item(32, 14).Value = -11651
$ws.Cells.Item(38, 8).Value = 106.125
$ws.Cells.Item(38, 9).Value = 106.125
$ws.Cells.Item(38, 11).Value = 318.375
$ws.Cells.Item(38, 13).Value = 53.625
$ws.Cells.Item(43, 8).Value = 12378.8
$ws.Cells.Item(43, 9).Value = 14882.167
$ws.Cells.Item(43, 11).Value = 14882.167
$ws.Cells.Item(43, 13).Value = -14813.167
$ws.Cells.Item(58, 8).Value = 3099
$ws.Cells.Item(58, 9).Value = 339
$ws.Cells.Item(58, 10).Value = 9999
$ws.Cells.Item(58, 11).Value = 1017
$ws.Cells.Item(58, 12).Value = 29997
$ws.Cells.Item(58, 13).Value = -867
$ws.Cells.Item(58, 14).Value = -30297
$ws.Cells.Item(64, 8).Value = 7632.6665
$ws.Cells.Item(64, 10).Value = 7499.25
$ws.Cells.Item(64, 12).Value = 7499.25
$ws.Cells.Item(64, 14).Value = -7995.25
$ws.Cells.Item(67, 8).Value = 7632.6665
$ws.Cells.Item(67, 10).Value = 7499.25
$ws.Cells.Item(67, 12).Value = 7499.25
$ws.Cells.Item(67, 14).Value = -9215.25
$ws.Cells.Item(69, 8).Value = 20000
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 20000
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 60000
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -61748
$ws.Cells.Item(70, 8).Value = 6742.4165
$ws.Cells.Item(70, 9).Value = 4974.5
$ws.Cells.Item(70, 10).Value = 7626.375
$ws.Cells.Item(70, 11).Value = 14923.5
$ws.Cells.Item(70, 12).Value = 22879.125
$ws.Cells.Item(70, 13).Value = -14653.5
$ws.Cells.Item(70, 14).Value = -23419.125
$ws.Cells.Item(72, 8).Value = 20000
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 20000
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 180000
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -188736
$ws.Cells.Item(73, 8).Value = 6742.4165
$ws.Cells.Item(73, 9).Value = 4974.5
$ws.Cells.Item(73, 10).Value = 7626.375
$ws.Cells.Item(73, 11).Value = 14923.5
$ws.Cells.Item(73, 12).Value = 22879.125
$ws.Cells.Item(73, 13).Value = -13987.5
$ws.Cells.Item(73, 14).Value = -24751.125
$ws.Cells.Item(80, 8).Value = 1464.0625
$ws.Cells.Item(80, 9).Value = 676.8570999999999
$ws.Cells.Item(80, 10).Value = 2076.3333
$ws.Cells.Item(80, 11).Value = 2030.5713
$ws.Cells.Item(80, 12).Value = 6228.999899999999
$ws.Cells.Item(80, 13).Value = -1032.5713
$ws.Cells.Item(80, 14).Value = -8224.999899999999
$ws.Cells.Item(83, 8).Value = 1464.0625
$ws.Cells.Item(83, 9).Value = 676.8570999999999
$ws.Cells.Item(83, 10).Value = 2076.3333
$ws.Cells.Item(83, 11).Value = 6091.7139
$ws.Cells.Item(83, 12).Value = 18686.9997
$ws.Cells.Item(83, 13).Value = -1099.7139
$ws.Cells.Item(83, 14).Value = -28670.9997
$ws.Cells.Item(96, 8).Value = 1470.25
$ws.Cells.Item(96, 9).Value = 954.6
$ws.Cells.Item(96, 10).Value = 2329.6667
$ws.Cells.Item(96, 11).Value = 2863.8
$ws.Cells.Item(96, 12).Value = 6989.000100000001
$ws.Cells.Item(96, 13).Value = -1490.8
$ws.Cells.Item(96, 14).Value = -9735.000100000001
$ws.Cells.Item(100, 8).Value = 2022.091
$ws.Cells.Item(100, 9).Value = 1884.7142
$ws.Cells.Item(100, 11).Value = 1884.7142
$ws.Cells.Item(100, 13).Value = -1343.7142
$ws.Cells.Item(102, 8).Value = 39933
$ws.Cells.Item(102, 10).Value = 39933
$ws.Cells.Item(102, 12).Value = 39933
$ws.Cells.Item(102, 14).Value = -46423
$ws.Cells.Item(106, 8).Value = 7399
$ws.Cells.Item(106, 9).Value = 7399
$ws.Cells.Item(106, 11).Value = 7399
$ws.Cells.Item(106, 13).Value = -6768
$ws.Cells.Item(116, 8).Value = 6511.9546
$ws.Cells.Item(116, 9).Value = 5890.769
$ws.Cells.Item(116, 11).Value = 5890.769
$ws.Cells.Item(116, 13).Value = -2448.769
$ws.Cells.Item(125, 8).Value = 961.5714
$ws.Cells.Item(125, 9).Value = 982.5
$ws.Cells.Item(125, 11).Value = 8842.5
$ws.Cells.Item(125, 13).Value = -6382.5
$ws.Cells.Item(132, 8).Value = 2868.6128
$ws.Cells.Item(132, 9).Value = 2890.6072
$ws.Cells.Item(132, 10).Value = 2663.3333
$ws.Cells.Item(132, 11).Value = 8671.821599999999
$ws.Cells.Item(132, 12).Value = 7989.999899999999
$ws.Cells.Item(132, 13).Value = -6141.821599999999
$ws.Cells.Item(132, 14).Value = -13049.9999
$ws.Cells.Item(135, 8).Value = 5440.154
$ws.Cells.Item(135, 9).Value = 1262.8572
$ws.Cells.Item(135, 11).Value = 11365.7148
$ws.Cells.Item(135, 13).Value = -8830.7148
$ws.Cells.Item(141, 8).Value = 4024.0908
$ws.Cells.Item(141, 9).Value = 3547
$ws.Cells.Item(141, 11).Value = 10641
$ws.Cells.Item(141, 13).Value = -5461
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 1000000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(45, 8).Value = 9366.767
$ws.Cells.Item(45, 9).Value = 12141.7
$ws.Cells.Item(45, 10).Value = 3816.9
$ws.Cells.Item(45, 11).Value = 12141.7
$ws.Cells.Item(45, 12).Value = 3816.9
$ws.Cells.Item(45, 13).Value = -11764.7
$ws.Cells.Item(45, 14).Value = -4570.9
$ws.Cells.Item(61, 8).Value = 9713.5
$ws.Cells.Item(61, 9).Value = 9572.200000000001
$ws.Cells.Item(61, 10).Value = 9949
$ws.Cells.Item(61, 11).Value = 9572.200000000001
$ws.Cells.Item(61, 12).Value = 9949
$ws.Cells.Item(61, 13).Value = -9360.200000000001
$ws.Cells.Item(61, 14).Value = -10373
$ws.Cells.Item(62, 8).Value = 28499
$ws.Cells.Item(62, 10).Value = 28499
$ws.Cells.Item(62, 12).Value = 28499
$ws.Cells.Item(62, 14).Value = -29747
$ws.Cells.Item(65, 8).Value = 28499
$ws.Cells.Item(65, 10).Value = 28499
$ws.Cells.Item(65, 12).Value = 85497
$ws.Cells.Item(65, 14).Value = -91737
$ws.Cells.Item(74, 8).Value = 1715.2916
$ws.Cells.Item(74, 10).Value = 1851
$ws.Cells.Item(74, 12).Value = 1851
$ws.Cells.Item(74, 14).Value = -3599
$ws.Cells.Item(77, 8).Value = 1715.2916
$ws.Cells.Item(77, 10).Value = 1851
$ws.Cells.Item(77, 12).Value = 9255
$ws.Cells.Item(77, 14).Value = -17991
$ws.Cells.Item(88, 8).Value = 2391.5
$ws.Cells.Item(88, 10).Value = 2588.6667
$ws.Cells.Item(88, 12).Value = 2588.6667
$ws.Cells.Item(88, 14).Value = -3400.6667
$ws.Cells.Item(91, 8).Value = 2391.5
$ws.Cells.Item(91, 10).Value = 2588.6667
$ws.Cells.Item(91, 12).Value = 2588.6667
$ws.Cells.Item(91, 14).Value = -5396.6667
$ws.Cells.Item(92, 8).Value = 55333.332
$ws.Cells.Item(92, 10).Value = 55333.332
$ws.Cells.Item(92, 12).Value = 55333.332
$ws.Cells.Item(92, 14).Value = -60325.332
$ws.Cells.Item(97, 8).Value = 1381.909
$ws.Cells.Item(97, 9).Value = 1222.3334
$ws.Cells.Item(97, 11).Value = 1222.3334
$ws.Cells.Item(97, 13).Value = -726.3334
$ws.Cells.Item(102, 8).Value = 3844.35
$ws.Cells.Item(102, 9).Value = 3666.0557
$ws.Cells.Item(102, 10).Value = 5449
$ws.Cells.Item(102, 11).Value = 3666.0557
$ws.Cells.Item(102, 12).Value = 5449
$ws.Cells.Item(102, 13).Value = -2044.0557
$ws.Cells.Item(102, 14).Value = -8693
$ws.Cells.Item(110, 8).Value = 2229
$ws.Cells.Item(110, 9).Value = 1852.9445
$ws.Cells.Item(110, 11).Value = 1852.9445
$ws.Cells.Item(110, 13).Value = 192.0554999999999
$ws.Cells.Item(132, 8).Value = 2845.7036
$ws.Cells.Item(132, 9).Value = 2450.4167
$ws.Cells.Item(132, 11).Value = 7351.250100000001
$ws.Cells.Item(132, 13).Value = -4821.250100000001
$ws.Cells.Item(136, 8).Value = 9713.5
$ws.Cells.Item(136, 9).Value = 9572.200000000001
$ws.Cells.Item(136, 10).Value = 9949
$ws.Cells.Item(136, 11).Value = 28716.6
$ws.Cells.Item(136, 12).Value = 29847
$ws.Cells.Item(136, 13).Value = -26166.6
$ws.Cells.Item(136, 14).Value = -34947
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1697.8334
$ws.Cells.Item(20, 9).Value = 1919.7
$ws.Cells.Item(20, 10).Value = 1254.1
$ws.Cells.Item(20, 11).Value = 1919.7
$ws.Cells.Item(20, 12).Value = 1254.1
$ws.Cells.Item(20, 13).Value = -1672.7
$ws.Cells.Item(20, 14).Value = -1748.1
$ws.Cells.Item(54, 8).Value = 5018.3335
$ws.Cells.Item(54, 9).Value = 5018.3335
$ws.Cells.Item(54, 11).Value = 5018.3335
$ws.Cells.Item(54, 13).Value = -4534.3335
$ws.Cells.Item(86, 8).Value = 5444
$ws.Cells.Item(86, 9).Value = 5107.643
$ws.Cells.Item(86, 11).Value = 5107.643
$ws.Cells.Item(86, 13).Value = -3984.643
$ws.Cells.Item(89, 8).Value = 5444
$ws.Cells.Item(89, 9).Value = 5107.643
$ws.Cells.Item(89, 11).Value = 25538.215
$ws.Cells.Item(89, 13).Value = -19922.215
$ws.Cells.Item(94, 8).Value = 909.5
$ws.Cells.Item(94, 9).Value = 949.3333
$ws.Cells.Item(94, 11).Value = 949.3333
$ws.Cells.Item(94, 13).Value = -498.3333
$ws.Cells.Item(100, 8).Value = 21425
$ws.Cells.Item(100, 10).Value = 21425
$ws.Cells.Item(100, 12).Value = 21425
$ws.Cells.Item(100, 14).Value = -23589
$ws.Cells.Item(105, 8).Value = 1488.5
$ws.Cells.Item(105, 9).Value = 1488.5
$ws.Cells.Item(105, 11).Value = 1488.5
$ws.Cells.Item(105, 13).Value = 258.5
$ws.Cells.Item(107, 8).Value = 2210.4517
$ws.Cells.Item(107, 9).Value = 2036.5
$ws.Cells.Item(107, 10).Value = 3834
$ws.Cells.Item(107, 11).Value = 2036.5
$ws.Cells.Item(107, 12).Value = 3834
$ws.Cells.Item(107, 13).Value = -116.5
$ws.Cells.Item(107, 14).Value = -7674
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 306.8889
$ws.Cells.Item(6, 9).Value = 306.8889
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 306.8889
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -193.8889
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(22, 8).Value = 1410.9
$ws.Cells.Item(22, 9).Value = 1410.9
$ws.Cells.Item(22, 11).Value = 1410.9
$ws.Cells.Item(22, 13).Value = -1060.9
$ws.Cells.Item(31, 8).Value = 3574.1555
$ws.Cells.Item(31, 9).Value = 3472.1562
$ws.Cells.Item(31, 11).Value = 3472.1562
$ws.Cells.Item(31, 13).Value = -3177.1562
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(34, 8).Value = 3574.1555
$ws.Cells.Item(34, 9).Value = 3472.1562
$ws.Cells.Item(34, 11).Value = 3472.1562
$ws.Cells.Item(34, 13).Value = -3270.1562
$ws.Cells.Item(50, 8).Value = 34940
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).ClearContents()
$ws.Cells.Item(64, 8).Value = 40000
$ws.Cells.Item(64, 10).Value = 40000
$ws.Cells.Item(64, 12).Value = 40000
$ws.Cells.Item(64, 14).Value = -40496
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).ClearContents()
$ws.Cells.Item(67, 8).Value = 40000
$ws.Cells.Item(67, 10).Value = 40000
$ws.Cells.Item(67, 12).Value = 40000
$ws.Cells.Item(67, 14).Value = -41716
$ws.Cells.Item(86, 8).Value = 30311420
$ws.Cells.Item(86, 9).Value = 83339980
$ws.Cells.Item(86, 11).Value = 83339980
$ws.Cells.Item(86, 13).Value = -83338857
$ws.Cells.Item(89, 8).Value = 30311420
$ws.Cells.Item(89, 9).Value = 83339980
$ws.Cells.Item(89, 11).Value = 416699900
$ws.Cells.Item(89, 13).Value = -416694284
$ws.Cells.Item(99, 8).Value = 4420.2856
$ws.Cells.Item(99, 9).Value = 3432.4443
$ws.Cells.Item(99, 11).Value = 3432.4443
$ws.Cells.Item(99, 13).Value = -1934.4443
$ws.Cells.Item(107, 8).Value = 1162.1111
$ws.Cells.Item(107, 9).Value = 994.2857
$ws.Cells.Item(107, 10).Value = 1749.5
$ws.Cells.Item(107, 11).Value = 994.2857
$ws.Cells.Item(107, 12).Value = 1749.5
$ws.Cells.Item(107, 13).Value = 925.7143
$ws.Cells.Item(107, 14).Value = -5589.5
$ws.Cells.Item(126, 8).Value = 4420.2856
$ws.Cells.Item(126, 9).Value = 3432.4443
$ws.Cells.Item(126, 11).Value = 10297.3329
$ws.Cells.Item(126, 13).Value = -7827.332900000001
$ws.Cells.Item(132, 8).Value = 3448
$ws.Cells.Item(132, 9).Value = 2904.5854
$ws.Cells.Item(132, 11).Value = 8713.7562
$ws.Cells.Item(132, 13).Value = -6183.7562
$ws.Cells.Item(134, 8).Value = 5867.1313
$ws.Cells.Item(134, 9).Value = 5375.5483
$ws.Cells.Item(134, 11).Value = 16126.6449
$ws.Cells.Item(134, 13).Value = -13591.6449
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3192.8823
$ws.Cells.Item(3, 9).Value = 2642.5
$ws.Cells.Item(3, 11).Value = 7927.5
$ws.Cells.Item(3, 13).Value = -7815.5
$ws.Cells.Item(5, 8).Value = 724.1667
$ws.Cells.Item(5, 9).Value = 669
$ws.Cells.Item(5, 11).Value = 2007
$ws.Cells.Item(5, 13).Value = -1895
$ws.Cells.Item(7, 8).Value = 690.8889
$ws.Cells.Item(7, 9).Value = 235.21739
$ws.Cells.Item(7, 10).Value = 3311
$ws.Cells.Item(7, 11).Value = 705.65217
$ws.Cells.Item(7, 12).Value = 9933
$ws.Cells.Item(7, 13).Value = -593.65217
$ws.Cells.Item(7, 14).Value = -10157
$ws.Cells.Item(12, 8).Value = 249.52777
$ws.Cells.Item(12, 9).Value = 209
$ws.Cells.Item(12, 11).Value = 627
$ws.Cells.Item(12, 13).Value = -454
$ws.Cells.Item(17, 8).Value = 1168.6923
$ws.Cells.Item(17, 9).Value = 354
$ws.Cells.Item(17, 11).Value = 1062
$ws.Cells.Item(17, 13).Value = -893
$ws.Cells.Item(23, 8).Value = 696.43475
$ws.Cells.Item(23, 9).Value = 517
$ws.Cells.Item(23, 11).Value = 1551
$ws.Cells.Item(23, 13).Value = -1316
$ws.Cells.Item(33, 8).Value = 95.454544
$ws.Cells.Item(33, 9).Value = 73.833336
$ws.Cells.Item(33, 10).Value = 121.4
$ws.Cells.Item(33, 11).Value = 443.000016
$ws.Cells.Item(33, 12).Value = 728.4000000000001
$ws.Cells.Item(33, 13).Value = -160.000016
$ws.Cells.Item(33, 14).Value = -1294.4
$ws.Cells.Item(34, 8).Value = 3158.9048
$ws.Cells.Item(34, 10).Value = 3306.75
$ws.Cells.Item(34, 12).Value = 9920.25
$ws.Cells.Item(34, 14).Value = -10088.25
$ws.Cells.Item(55, 8).Value = 5661
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(68, 8).Value = 398.33334
$ws.Cells.Item(68, 9).Value = 418.6
$ws.Cells.Item(68, 11).Value = 1255.8
$ws.Cells.Item(68, 13).Value = -444.8000000000002
$ws.Cells.Item(71, 8).Value = 398.33334
$ws.Cells.Item(71, 9).Value = 418.6
$ws.Cells.Item(71, 11).Value = 3767.4
$ws.Cells.Item(71, 13).Value = 288.5999999999999
$ws.Cells.Item(75, 8).Value = 2324.3
$ws.Cells.Item(75, 10).Value = 2235.7144
$ws.Cells.Item(75, 12).Value = 6707.1432
$ws.Cells.Item(75, 14).Value = -8703.143199999999
$ws.Cells.Item(78, 8).Value = 2324.3
$ws.Cells.Item(78, 10).Value = 2235.7144
$ws.Cells.Item(78, 12).Value = 20121.4296
$ws.Cells.Item(78, 14).Value = -30105.4296
$ws.Cells.Item(92, 8).Value = 311.4
$ws.Cells.Item(92, 9).Value = 282.33334
$ws.Cells.Item(92, 10).Value = 323.85715
$ws.Cells.Item(92, 11).Value = 847.0000200000001
$ws.Cells.Item(92, 12).Value = 971.5714499999999
$ws.Cells.Item(92, 13).Value = 400.9999799999999
$ws.Cells.Item(92, 14).Value = -3467.57145
$ws.Cells.Item(96, 8).Value = 0
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 11).Value = 0
$ws.Cells.Item(96, 13).ClearContents()
$ws.Cells.Item(98, 8).Value = 1993.8889
$ws.Cells.Item(98, 10).Value = 2634.2
$ws.Cells.Item(98, 12).Value = 7902.599999999999
$ws.Cells.Item(98, 14).Value = -10898.6
$ws.Cells.Item(114, 8).Value = 1957.5
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 1957.5
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 5872.5
$ws.Cells.Item(114, 13).ClearContents()
$ws.Cells.Item(114, 14).Value = -12380.5
$ws.Cells.Item(117, 8).Value = 2772.818
$ws.Cells.Item(117, 10).Value = 3052.4443
$ws.Cells.Item(117, 12).Value = 9157.332900000001
$ws.Cells.Item(117, 14).Value = -16041.3329
$ws.Cells.Item(121, 8).Value = 1189.091
$ws.Cells.Item(121, 10).Value = 1208.2
$ws.Cells.Item(121, 12).Value = 3624.6
$ws.Cells.Item(121, 14).Value = -6244.6
$ws.Cells.Item(129, 8).Value = 1301
$ws.Cells.Item(129, 9).Value = 969.3333
$ws.Cells.Item(129, 11).Value = 2907.9999
$ws.Cells.Item(129, 13).Value = 2092.0001
$ws.Cells.Item(132, 8).Value = 1235.6364
$ws.Cells.Item(132, 10).Value = 1939.8
$ws.Cells.Item(132, 12).Value = 17458.2
$ws.Cells.Item(132, 14).Value = -22518.2
$ws.Cells.Item(133, 8).Value = 15999
$ws.Cells.Item(133, 9).Value = 11498.5
$ws.Cells.Item(133, 11).Value = 34495.5
$ws.Cells.Item(133, 13).Value = -29435.5
$ws.Cells.Item(135, 8).Value = 724.1667
$ws.Cells.Item(135, 9).Value = 669
$ws.Cells.Item(135, 11).Value = 6021
$ws.Cells.Item(135, 13).Value = -3486
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3428.2222
$ws.Cells.Item(80, 10).Value = 4326
$ws.Cells.Item(80, 12).Value = 4326
$ws.Cells.Item(80, 14).Value = -6322
$ws.Cells.Item(83, 8).Value = 3428.2222
$ws.Cells.Item(83, 10).Value = 4326
$ws.Cells.Item(83, 12).Value = 21630
$ws.Cells.Item(83, 14).Value = -31614
$ws.Cells.Item(97, 8).Value = 1150.9166
$ws.Cells.Item(97, 9).Value = 1185
$ws.Cells.Item(97, 11).Value = 1185
$ws.Cells.Item(97, 13).Value = -689
$ws.Cells.Item(102, 8).Value = 1862.625
$ws.Cells.Item(102, 9).Value = 1771.5714
$ws.Cells.Item(102, 10).Value = 2500
$ws.Cells.Item(102, 11).Value = 1771.5714
$ws.Cells.Item(102, 12).Value = 2500
$ws.Cells.Item(102, 13).Value = -149.5714
$ws.Cells.Item(102, 14).Value = -5744
$ws.Cells.Item(126, 8).Value = 4410.037
$ws.Cells.Item(126, 9).Value = 4109.067
$ws.Cells.Item(126, 10).Value = 4786.25
$ws.Cells.Item(126, 11).Value = 12327.201
$ws.Cells.Item(126, 12).Value = 14358.75
$ws.Cells.Item(126, 13).Value = -9857.201000000001
$ws.Cells.Item(126, 14).Value = -19298.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 11194.353
$ws.Cells.Item(7, 10).Value = 4886.4287
$ws.Cells.Item(7, 12).Value = 4886.4287
$ws.Cells.Item(7, 14).Value = -5110.4287
$ws.Cells.Item(9, 8).Value = 12749.5
$ws.Cells.Item(9, 9).Value = 6999
$ws.Cells.Item(9, 10).Value = 14666.333
$ws.Cells.Item(9, 11).Value = 6999
$ws.Cells.Item(9, 12).Value = 14666.333
$ws.Cells.Item(9, 13).Value = -6775
$ws.Cells.Item(9, 14).Value = -15114.333
$ws.Cells.Item(16, 8).Value = 1484.125
$ws.Cells.Item(16, 10).Value = 2692.5
$ws.Cells.Item(16, 12).Value = 2692.5
$ws.Cells.Item(16, 14).Value = -3032.5
$ws.Cells.Item(22, 8).Value = 2859.4546
$ws.Cells.Item(22, 9).Value = 2793.65
$ws.Cells.Item(22, 10).Value = 2960.6924
$ws.Cells.Item(22, 11).Value = 2793.65
$ws.Cells.Item(22, 12).Value = 2960.6924
$ws.Cells.Item(22, 13).Value = -2498.65
$ws.Cells.Item(22, 14).Value = -3550.6924
$ws.Cells.Item(27, 8).Value = 2859.4546
$ws.Cells.Item(27, 9).Value = 2793.65
$ws.Cells.Item(27, 10).Value = 2960.6924
$ws.Cells.Item(27, 11).Value = 2793.65
$ws.Cells.Item(27, 12).Value = 2960.6924
$ws.Cells.Item(27, 13).Value = -2686.65
$ws.Cells.Item(27, 14).Value = -3174.6924
$ws.Cells.Item(40, 8).Value = 16698.5
$ws.Cells.Item(40, 9).Value = 13382.857
$ws.Cells.Item(40, 11).Value = 13382.857
$ws.Cells.Item(40, 13).Value = -13246.857
$ws.Cells.Item(46, 8).Value = 3806.7
$ws.Cells.Item(46, 9).Value = 3882.4443
$ws.Cells.Item(46, 11).Value = 3882.4443
$ws.Cells.Item(46, 13).Value = -3694.4443
$ws.Cells.Item(55, 8).Value = 224.59259
$ws.Cells.Item(55, 10).Value = 167
$ws.Cells.Item(55, 12).Value = 167
$ws.Cells.Item(55, 14).Value = -513
$ws.Cells.Item(61, 8).Value = 1585.5
$ws.Cells.Item(61, 9).Value = 1760.625
$ws.Cells.Item(61, 10).Value = 1352
$ws.Cells.Item(61, 11).Value = 1760.625
$ws.Cells.Item(61, 12).Value = 1352
$ws.Cells.Item(61, 13).Value = -1558.625
$ws.Cells.Item(61, 14).Value = -1756
$ws.Cells.Item(93, 8).Value = 1573.2727
$ws.Cells.Item(93, 9).Value = 1227.7142
$ws.Cells.Item(93, 10).Value = 2178
$ws.Cells.Item(93, 11).Value = 1227.7142
$ws.Cells.Item(93, 12).Value = 2178
$ws.Cells.Item(93, 13).Value = 20.28580000000011
$ws.Cells.Item(93, 14).Value = -4674
$ws.Cells.Item(100, 8).Value = 1992.8
$ws.Cells.Item(100, 9).Value = 1992.8
$ws.Cells.Item(100, 11).Value = 1992.8
$ws.Cells.Item(100, 13).Value = -1451.8
$ws.Cells.Item(113, 8).Value = 1585.5
$ws.Cells.Item(113, 9).Value = 1760.625
$ws.Cells.Item(113, 10).Value = 1352
$ws.Cells.Item(113, 11).Value = 1760.625
$ws.Cells.Item(113, 12).Value = 1352
$ws.Cells.Item(113, 13).Value = 409.375
$ws.Cells.Item(113, 14).Value = -5692
$ws.Cells.Item(115, 8).Value = 89500
$ws.Cells.Item(115, 10).Value = 89500
$ws.Cells.Item(115, 12).Value = 89500
$ws.Cells.Item(115, 14).Value = -91850
$ws.Cells.Item(126, 8).Value = 11194.353
$ws.Cells.Item(126, 10).Value = 4886.4287
$ws.Cells.Item(126, 12).Value = 14659.2861
$ws.Cells.Item(126, 14).Value = -19599.2861
$ws.Cells.Item(132, 8).Value = 3436.9375
$ws.Cells.Item(132, 9).Value = 2744.9092
$ws.Cells.Item(132, 10).Value = 4959.4
$ws.Cells.Item(132, 11).Value = 8234.7276
$ws.Cells.Item(132, 12).Value = 14878.2
$ws.Cells.Item(132, 13).Value = -5704.7276
$ws.Cells.Item(132, 14).Value = -19938.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 262248.38
$ws.Cells.Item(4, 9).Value = 272997
$ws.Cells.Item(4, 10).Value = 251499.75
$ws.Cells.Item(4, 11).Value = 272997
$ws.Cells.Item(4, 12).Value = 251499.75
$ws.Cells.Item(4, 13).Value = -272884
$ws.Cells.Item(4, 14).Value = -251725.75
$ws.Cells.Item(41, 8).Value = 10499.5
$ws.Cells.Item(41, 10).Value = 10499.5
$ws.Cells.Item(41, 12).Value = 10499.5
$ws.Cells.Item(41, 14).Value = -11279.5
$ws.Cells.Item(46, 8).Value = 58895.332
$ws.Cells.Item(46, 9).Value = 48888
$ws.Cells.Item(46, 10).Value = 63899
$ws.Cells.Item(46, 11).Value = 48888
$ws.Cells.Item(46, 12).Value = 63899
$ws.Cells.Item(46, 13).Value = -48657
$ws.Cells.Item(46, 14).Value = -64361
$ws.Cells.Item(62, 8).Value = 27999.6
$ws.Cells.Item(62, 10).Value = 37500
$ws.Cells.Item(62, 12).Value = 37500
$ws.Cells.Item(62, 14).Value = -38748
$ws.Cells.Item(65, 8).Value = 27999.6
$ws.Cells.Item(65, 10).Value = 37500
$ws.Cells.Item(65, 12).Value = 187500
$ws.Cells.Item(65, 14).Value = -193740
$ws.Cells.Item(81, 8).Value = 3137.2666
$ws.Cells.Item(81, 9).Value = 3579
$ws.Cells.Item(81, 11).Value = 7158
$ws.Cells.Item(81, 13).Value = -6097
$ws.Cells.Item(84, 8).Value = 3137.2666
$ws.Cells.Item(84, 9).Value = 3579
$ws.Cells.Item(84, 11).Value = 35790
$ws.Cells.Item(84, 13).Value = -30486
$ws.Cells.Item(100, 8).Value = 1643.9445
$ws.Cells.Item(100, 10).Value = 1616.8572
$ws.Cells.Item(100, 12).Value = 3233.7144
$ws.Cells.Item(100, 14).Value = -4315.7144
$ws.Cells.Item(107, 8).Value = 3450.3333
$ws.Cells.Item(107, 9).Value = 4056.6365
$ws.Cells.Item(107, 11).Value = 12169.9095
$ws.Cells.Item(107, 13).Value = -10249.9095
$ws.Cells.Item(126, 8).Value = 2452.9355
$ws.Cells.Item(126, 9).Value = 1814.5294
$ws.Cells.Item(126, 11).Value = 5443.5882
$ws.Cells.Item(126, 13).Value = -2973.5882
$ws.Cells.Item(127, 8).Value = 83331.664
$ws.Cells.Item(127, 9).Value = 98997.5
$ws.Cells.Item(127, 10).Value = 52000
$ws.Cells.Item(127, 11).Value = 98997.5
$ws.Cells.Item(127, 12).Value = 52000
$ws.Cells.Item(127, 13).Value = -94037.5
$ws.Cells.Item(127, 14).Value = -61920
$ws.Cells.Item(132, 8).Value = 2332.475
$ws.Cells.Item(132, 9).Value = 2574
$ws.Cells.Item(132, 10).Value = 1883.9286
$ws.Cells.Item(132, 11).Value = 7722
$ws.Cells.Item(132, 12).Value = 5651.7858
$ws.Cells.Item(132, 13).Value = -5192
$ws.Cells.Item(132, 14).Value = -10711.7858
$ws.Cells.Item(134, 8).Value = 58895.332
$ws.Cells.Item(134, 9).Value = 48888
$ws.Cells.Item(134, 10).Value = 63899
$ws.Cells.Item(134, 11).Value = 146664
$ws.Cells.Item(134, 12).Value = 191697
$ws.Cells.Item(134, 13).Value = -144129
$ws.Cells.Item(134, 14).Value = -196767